$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("J27").Value = "total distance"
$ws2.Range("J29").Value = 1
$ws2.Range("K29").Value = 2
$ws2.Range("L29").Value = 3

$ws1 = $wb.Worksheets.Item("Sheet1")
$co = $ws1.ChartObjects(2)
$ch = $co.Chart
$sc = $ch.SeriesCollection()
$ser1 = $sc.Item(1)
$ser1.Formula = "=SERIES(Sheet2!`$I`$29,Sheet2!`$J`$27:`$L`$27,Sheet2!`$J`$29:`$L`$29,1)"
Write-Host ("After: " + $ser1.Formula)
